# Generate Report for Handback
# Refresh the handback timestamps for the first tracked file
# (2f80e1f8-af23-4406-a0df-9a04aefade6a) across the Overview, zh-cn and
# de-de sheets to reflect a newly-generated handback report.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: Latest HO Xliff Generate Date for the first file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-14 03:07:24"

# --- zh-cn sheet: Correspond Handoff/Handback Datetime for the first file ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-14 03:07:16"
$wsZhCn.Range("K2").Value = "2016-08-14 03:07:45"

# --- de-de sheet: Correspond Handoff/Handback Datetime for the first file ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-14 03:07:24"
$wsDeDe.Range("K2").Value = "2016-08-14 03:07:56"
